$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-32 (column C): fitness value 7590 -> 7598
$ws.Range("C2:C32").Value = 7598

# Rows 147-220 (column C): fitness value 7573 -> 7590
$ws.Range("C147:C220").Value = 7590
